# Adds table notes to the "Enemies" sheet (renames Type/Damage Type columns,
# reworks the Rarity header and adds a Sprite note), and introduces two new
# reference sheets: "Damage Types" and "Attack Type".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Enemies sheet: update the Rarity header note text first (J1) -- this is
#    the first edit chronologically so new shared strings are appended in
#    the same order the original author produced them in.
# ---------------------------------------------------------------------------
$enemies = $wb.Worksheets.Item("Enemies")
$enemies.Range("J1").Value = "Rarity (0-100)"

# ---------------------------------------------------------------------------
# 2. Add the "Damage Types" sheet right after "Enemies" and seed its id
#    header.
# ---------------------------------------------------------------------------
$afterEnemies = $wb.Worksheets.Item($wb.Worksheets.Count)
$damageTypes = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterEnemies)
$damageTypes.Name = "Damage Types"
$damageTypes.Range("A1").Value = "Id"
$damageTypes.Range("B1").Value = "Name"

# ---------------------------------------------------------------------------
# 3. Enemies sheet: swap the old "Type" column header/note for "Sprite".
# ---------------------------------------------------------------------------
$enemies.Range("C1").Value = "Sprite"
$enemies.Range("C2").Value = "determines which sprite is rendered for an enemy"

# ---------------------------------------------------------------------------
# 4. Populate "Damage Types" data rows (B = damage type name, D = the
#    attack type it belongs to).
# ---------------------------------------------------------------------------
$damageTypes.Range("B2").Value = "Melee"
$damageTypes.Range("B3").Value = "Projectile"
$damageTypes.Range("B4").Value = "Elemental-Fire"
$damageTypes.Range("B5").Value = "Elemental-Ice"
$damageTypes.Range("B6").Value = "Elemental-Electric"
$damageTypes.Range("B7").Value = "Elemental-Fire"
$damageTypes.Range("B8").Value = "Elemental-Ice"
$damageTypes.Range("B9").Value = "Elemental-Electric"

$damageTypes.Range("D2").Value = "Melee"
$damageTypes.Range("D3").Value = "Projectile"
$damageTypes.Range("D4").Value = "Melee"
$damageTypes.Range("D5").Value = "Melee"
$damageTypes.Range("D6").Value = "Melee"
$damageTypes.Range("D7").Value = "Projectile"
$damageTypes.Range("D8").Value = "Projectile"
$damageTypes.Range("D9").Value = "Projectile"

# Headers for the two attack-type lookup columns (typed last).
$damageTypes.Range("C1").Value = "Attack_Type_Id"
$damageTypes.Range("D1").Value = "Attack_Type_Name"

# ---------------------------------------------------------------------------
# 5. Enemies sheet: the old "Damage Type" column now stores a damage_type_id.
# ---------------------------------------------------------------------------
$enemies.Range("D1").Value = "damage_type_id"

# ---------------------------------------------------------------------------
# 6. Add the "Attack Type" sheet after "Damage Types".
# ---------------------------------------------------------------------------
$afterDamageTypes = $wb.Worksheets.Item($wb.Worksheets.Count)
$attackType = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterDamageTypes)
$attackType.Name = "Attack Type"
$attackType.Range("A1").Value = "id"
$attackType.Range("B1").Value = "Name"
$attackType.Range("B2").Value = "Melee"
$attackType.Range("B3").Value = "Projectile"

# ---------------------------------------------------------------------------
# 7. Column widths on "Damage Types" (best-effort match of the author's
#    manual column resize / autofit).
# ---------------------------------------------------------------------------
$damageTypes.Columns.Item(2).ColumnWidth = 16.59
$damageTypes.Columns.Item(3).ColumnWidth = 13.88
$damageTypes.Columns.Item(4).ColumnWidth = 16.59

# ---------------------------------------------------------------------------
# 8. Restore per-sheet selections to match the author's final cursor
#    positions, then activate "Damage Types" as the active tab.
# ---------------------------------------------------------------------------
$enemies.Activate()
$enemies.Range("D2").Select()

$attackType.Activate()
$attackType.Range("F39").Select()

$damageTypes.Activate()
$damageTypes.Range("V15").Select()

Write-Host "Applied Damage Types / Attack Type sheet changes"
